# Presentation almost ready, Patrones by @mgtapia missing
# The "Patrones" slide (last slide, slide8.xml / sldId 263) is removed.
$p = $ppt.ActivePresentation

$lastIndex = $p.Slides.Count
$s = $p.Slides.Item($lastIndex)
$s.Delete()
